# Weekly update: insert a new data row at row 53 (pushing existing rows
# 53-82 down to 54-83) and populate it with this week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 53, shifting everything below
# (rows 53-82) down by one (to 54-83). This also grows the used range
# from A1:T82 to A1:T83.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new weekly record.
$ws.Range("A53").Value = 11
$ws.Range("B53").Value = "Vega Monumental Concepción"
$ws.Range("C53").Value = "Bíobío"
$ws.Range("D53").Value = 44582
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100101
$ws.Range("H53").Value = "Berries"
$ws.Range("I53").Value = 100101001
$ws.Range("J53").Value = "Arándano (blue)"
$ws.Range("K53").Value = "Sin especificar"
$ws.Range("L53").Value = "Primera"
$ws.Range("M53").Value = 250
$ws.Range("N53").Value = 2800
$ws.Range("O53").Value = 3000
$ws.Range("P53").Value = 2920
$ws.Range("Q53").Value = "$/bandeja 2 kilos"
$ws.Range("R53").Value = "Provincia de Linares"
$ws.Range("S53").Value = 1460
$ws.Range("T53").Value = 2
